$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
  D2 = '65.419.86'
  E2 = '  -0.14%  '
  D3 = '3.561.72'
  E3 = '  +3.61%  '
  E4 = '  -0.04%  '
  D5 = '598.14'
  E5 = '  +0.73%  '
  D6 = '140.53'
  E6 = '  +3.73%  '
  D7 = '3.563.22'
  E7 = '  +3.68%  '
  E8 = '  +0.10%  '
  E9 = '  +1.25%  '
  E10 = '  +3.44%  '
  D11 = '7.17'
  E11 = '  -5.16%  '
  E12 = '  +3.43%  '
  D13 = '4.160.92'
  E13 = '  +3.56%  '
  B15 = 'WrappedEther'
  C15 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
  D15 = '3.557.70'
  E15 = '  +3.02%  '
  B16 = 'Avalanche'
  C16 = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
  D16 = '27.03'
  E16 = '  +2.12%  '
  D18 = '65.340.37'
  E18 = '  -0.14%  '
  E19 = '  +4.29%  '
  D20 = '5.86'
  E20 = '  +1.79%  '
  D21 = '14.21'
  E21 = '  +3.59%  '
  D22 = '397.16'
  E22 = '  +1.30%  '
  E23 = '  +4.63%  '
  B24 = 'Litecoin'
  C24 = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
  D24 = '74.74'
  E24 = '  +2.28%  '
  B25 = 'WrappedeETH'
  C25 = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
  D25 = '3.700.57'
  E25 = '  +3.39%  '
  E26 = '  -0.09%  '
  E27 = '  +9.26%  '
  D28 = '7.76'
  E28 = '  +6.11%  '
  D29 = '0.998'
  E30 = '  +0.55%  '
  D31 = '8.27'
  E31 = '  +1.53%  '
  D32 = '3.571.73'
  E32 = '  +3.70%  '
  D33 = '24.09'
  E33 = '  +6.12%  '
  E35 = '  +1.69%  '
  E36 = '  +1.64%  '
  E37 = '  +2.90%  '
  D38 = '168.61'
  E38 = '  -1.95%  '
  E39 = '  +1.76%  '
  E40 = '  +4.00%  '
  D41 = '0.0804'
  E41 = '  +3.95%  '
  D42 = '0.829'
  E42 = '  +1.87%  '
  D43 = '26.39'
  E43 = '  +15.44%  '
  D44 = '43.00'
  E44 = '  -1.20%  '
  E45 = '  +0.12%  '
  D46 = '4.44'
  E46 = '  +0.67%  '
  E47 = '  +4.34%  '
  E48 = '  +8.70%  '
  D49 = '6.82'
  E49 = '  +4.54%  '
  D50 = '2.405.36'
  E50 = '  +9.83%  '
  D51 = '2.12'
  E51 = '  +0.18%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
